$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("B10").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C10").Value = "3577649 - Carlos Angelo Nunes"

# Row 13
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "01/01/2013"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "01/01/2013"
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C15").Value = "3577649 - Carlos Angelo Nunes"
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Range("A16").Value = "Syllabus:"
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Rows.Item(17).AutoFit()

# Row 18
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
$ws.Range("C19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Critério`nMF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio."
$ws.Range("C20").Value = "Critério`nMF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio."

# Row 21
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Não será oferecida recuperação."
$ws.Range("C21").Value = "Não será oferecida recuperação."
$ws.Rows.Item(21).RowHeight = 120

# Row 22
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Rows.Item(22).AutoFit()

# Row 23
$ws.Range("A23").Clear()
$ws.Range("B23").Value = "LOM3070 -  Estágio Supervisionado  (Requisito)`n"
$ws.Range("C23").Value = "LOM3070 -  Estágio Supervisionado  (Requisito)`n"
$ws.Rows.Item(23).RowHeight = 30

# Fix column styles for brand-new B/C cells (engine col-style lookup bug: overlapping <col> ranges)
$ws.Range("B3").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fix styles for date-like text cells munged by NumberFormat-as-text change
$ws.Range("B3").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Remove now-obsolete trailing rows (old 24 & 25, content folded into rows 13-23 above)
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(24).Delete()

Write-Host "LOM3097 sheet restructuring applied"